# Delete the "Project Options" slide (originally the 2nd slide, sldId 650).
# This slide told students to "Remove this slide in your submission", so it
# is being removed from the deck.

$p = $ppt.ActivePresentation

$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $title = ""
    try {
        $title = $slide.Shapes.Item(1).TextFrame.TextRange.Text
    } catch {
        $title = ""
    }
    if ($title -eq "Project Options") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    # Fallback: the slide was originally the second slide in the deck.
    $targetIndex = 2
}

$p.Slides.Item($targetIndex).Delete()
